# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51.
# Numeric-looking price strings get NumberFormat "@" first so Excel keeps them as
# text (matching the source inlineStr cells) instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.779.33'
$ws.Range("E2").Value = '  -0.74%  '
$ws.Range("D3").Value = '1.626.05'
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.18'
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5113'
$ws.Range("E6").Value = '  +0.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2566'
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06325'
$ws.Range("E9").Value = '  -0.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.42'
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07781'
$ws.Range("E11").Value = '  +0.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.239'
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("E13").Value = '  -0.50%  '
$ws.Range("D14").Value = '1.848.89'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5527'
$ws.Range("E15").Value = '  +1.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.51'
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("D17").Value = '0.0₅7477'
$ws.Range("E17").Value = '  -2.76%  '
$ws.Range("D18").Value = '25.795.64'
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.419'
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '194.20'
$ws.Range("E21").Value = '  -2.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.767'
$ws.Range("E22").Value = '  -1.23%  '
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.862'
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.29'
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1248'
$ws.Range("E27").Value = '  +4.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.50'
$ws.Range("E28").Value = '  -0.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.705'
$ws.Range("E29").Value = '  -1.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.240'
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04859'
$ws.Range("E31").Value = '  -0.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.238'
$ws.Range("E32").Value = '  -0.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.170'
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.539'
$ws.Range("E34").Value = '  +0.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.361'
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8946'
$ws.Range("E36").Value = '  -1.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5515'
$ws.Range("E37").Value = '  +1.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.541'
$ws.Range("E38").Value = '  -1.51%  '
$ws.Range("D39").Value = '1.113.52'
$ws.Range("E39").Value = '  -1.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01547'
$ws.Range("E40").Value = '  -0.95%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.528'
$ws.Range("E42").Value = '  +2.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7958'
$ws.Range("E43").Value = '  -1.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.25'
$ws.Range("E44").Value = '  -1.73%  '
$ws.Range("D45").Value = '1.774.09'
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = '0.0₈115'
$ws.Range("E46").Value = '  -8.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4425'
$ws.Range("E47").Value = '  -2.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("E48").Value = '  -0.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.58'
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.571'
$ws.Range("E51").Value = '  +3.34%  '